$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 754, shifting the existing rows (754-795) down to (755-796)
$ws.Rows.Item(754).Insert()

# Populate the newly inserted row with the new data point for 2026/02/04 16:00
# Use text number format first so the date-looking string isn't auto-converted
# to a date serial, then clear the formatting so no style index is left behind
# (matches the rest of the data rows, which carry no explicit style).
$ws.Range("A754").NumberFormat = "@"
$ws.Range("A754").Value = "2026/02/04"
$ws.Range("A754").ClearFormats()

$ws.Range("B754").Value = "水"
$ws.Range("C754").Value = 16
$ws.Range("D754").Value = 35

"done"
